$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1 (summary paragraph): insert "computer " so the sentence reads
# "...a diverse set of computer software & languages." instead of
# "...a diverse set of software & languages."
# ---------------------------------------------------------------------
$summary = $d.Content
$summary.Find.Execute("a diverse set of software & languages.", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "a diverse set of computer software & languages.", 1)

# ---------------------------------------------------------------------
# Edit 2 (experience dates): the "Quality Assurance Analyst" entry
# ("May 2021 - Present") now has an end date, so its "Present" becomes
# "February 2022". The current role ("February 2022 - Present") must
# stay untouched, so scope the search to start right after "May 2021".
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("May 2021")
$scope = $d.Range($anchor.End, $d.Content.End)
$scope.Find.Execute("Present", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "February 2022", 1)
